$wb = $excel.ActiveWorkbook

$wsConversion = $wb.Worksheets.Item("conversionAssets")
$wsStorage = $wb.Worksheets.Item("storageAssets")

# --- conversionAssets (sheet3): new "ambientTempType" column (I) ---
# Enter AIR / GROUND before the header so the shared-string table fills
# up in the same order as the source workbook (AIR, GROUND, ambientTempType).
$wsConversion.Range("I2").Value = "AIR"
$wsConversion.Range("I3").Value = "AIR"
$wsConversion.Range("I4").Value = "AIR"
$wsConversion.Range("I6").Value = "AIR"
$wsConversion.Range("I9").Value = "GROUND"
$wsConversion.Range("I1").Value = "ambientTempType"

# capacityElectricity_kW for DH_heat_pump_HT_S bumped 100 -> 200
$wsConversion.Range("E9").Value = 200

# --- storageAssets (sheet4): new "ambientTempType" column (N) ---
$wsStorage.Range("N3").Value = "AIR"
$wsStorage.Range("N4").Value = "AIR"
$wsStorage.Range("N5").Value = "AIR"
$wsStorage.Range("N6").Value = "AIR"
$wsStorage.Range("N7").Value = "AIR"
$wsStorage.Range("N8").Value = "AIR"
$wsStorage.Range("N9").Value = "AIR"
$wsStorage.Range("N11").Value = "AIR"
$wsStorage.Range("N12").Value = "GROUND"
$wsStorage.Range("N13").Value = "GROUND"
$wsStorage.Range("N1").Value = "ambientTempType"

# widen the new column a touch, like the source file
$wsStorage.Columns.Item(14).ColumnWidth = 11.79

# --- sheet/view bookkeeping: conversionAssets becomes the active tab ---
$wsStorage.Range("H27").Select()
$wsConversion.Activate()
$wsConversion.Range("G27").Select()
